$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (D2: Target cluster -> ECs; E2 -> 3; G2,H2 updated; L2..T2 updated)
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.520808
$ws.Range("H2").Value = 7.562424
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1278803333333333
$ws.Range("N2").Value = 0.383641
$ws.Range("O2").Value = 0.009974564977605908
$ws.Range("P2").Value = 0.009974564977605908
$ws.Range("Q2").Value = 0.3223617673093334
$ws.Range("R2").Value = 2.901255905784
$ws.Range("S2").Value = 0.009974564977605908
$ws.Range("T2").Value = 0.009974564977605908

# Update row 3 (D3: Target cluster -> FAPs; E3 -> 3; G3,H3 updated; K3 -> 3; M3..T3 updated)
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.520808
$ws.Range("H3").Value = 7.562424
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.61985133333333
$ws.Range("N3").Value = 31.859554
$ws.Range("O3").Value = 0.8283400145723324
$ws.Range("P3").Value = 0.8283400145723324
$ws.Range("Q3").Value = 26.77060619987734
$ws.Range("R3").Value = 240.935455798896
$ws.Range("S3").Value = 0.8283400145723324
$ws.Range("T3").Value = 0.8283400145723324

# Add new row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt2"
$ws.Range("C4").Value = "Fzd2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.520808
$ws.Range("H4").Value = 7.562424
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.072911
$ws.Range("N4").Value = 6.218733
$ws.Range("O4").Value = 0.1616854204500617
$ws.Range("P4").Value = 0.1616854204500617
$ws.Range("Q4").Value = 5.225410632088
$ws.Range("R4").Value = 47.02869568879201
$ws.Range("S4").Value = 0.1616854204500617
$ws.Range("T4").Value = 0.1616854204500617
